$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J2:J11 placed in J12, bold font (11pt)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 11

# Row 14-17: summary labels + stats, bold 12pt font, vertical centered
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$statRange = $ws.Range("B14:B17")
$statRange.Font.Bold = $true
$statRange.Font.Size = 12
$statRange.VerticalAlignment = -4108

# Rows 14-17 end up taller because of the larger 12pt bold font (matches ht="15.6" in target)
$ws.Range("A14:B17").RowHeight = 15.6

# Set selection to match target (activeCell J12)
$ws.Range("J12").Select() | Out-Null

# Page setup (portrait, paper size 9 = A4) as in the target workbook
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9
